# "9th Stab - Cosmetic Changes"
#
# This workbook rolls a weekly analyst-rating snapshot: column B holds the
# newest period ("Jun_13" before this edit) and the rest of the sheet mostly
# holds the placeholder "UN" per broker, except for a single broker
# (row 7, Morgan Stanley) that has a real rating note highlighted in B7.
#
# The edit rolls the window forward by inserting two brand-new snapshot
# columns ("Jun_17" and "Jun_15") to the left of the existing data, pushing
# the old "Jun_13" column from B to D and the old "UN" header column from C
# to E. The new columns are seeded with the same "UN" placeholder used
# everywhere else, except Morgan Stanley's row gets a fresh note in the new
# "Jun_15" column, and the highlighted note itself stays with the shifted
# "Jun_13" column (now D) instead of staying in B.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: insert two new columns before the old column C, each seeded as a
# copy of column B (this duplicates the "UN" placeholder pattern already used
# throughout column B into the two new columns, and shifts the old column C
# two slots to the right, landing on E).
$ws.Columns.Item(2).Copy()
$ws.Columns.Item(3).Insert()
$ws.Columns.Item(2).Copy()
$ws.Columns.Item(3).Insert()

# Step 2: label the new snapshot columns.
$ws.Range("C1").Value = "Jun_15"
$ws.Range("B1").Value = "Jun_17"

# Step 3: Morgan Stanley (row 7) gets a brand new rating note in the newly
# inserted "Jun_15" column (C7) ...
$ws.Range("C7").Value = "6/13/2018,Reiterates,Neutral,$93.00"

# ... while the newest "Jun_17" column (B7) just falls back to the generic
# "UN" placeholder used by every other broker, and the old highlighted note
# (together with its highlight fill) stays attached to the old "Jun_13"
# column which is now D7 - so B7/C7 need to lose the highlight that the
# column copy above carried over from the original B7.
$ws.Range("B7").Value = "UN"
$ws.Range("B7").ClearFormats()
$ws.Range("C7").ClearFormats()

# Step 4: the old column C carried an explicit custom width (8.0 units).
# Re-apply that same width to C, D and E so all three keep explicit column
# formatting after the insert (matching Excel's own behaviour of carrying
# column formatting across an insert).
$ws.Columns.Item(3).ColumnWidth = 7.166666666666667
$ws.Columns.Item(4).ColumnWidth = 7.166666666666667
$ws.Columns.Item(5).ColumnWidth = 7.166666666666667
